$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "y (t/cu.m.)" unit-weight values in column E for rows 4-19
$values = 1.8, 1.6, 1.7, 1.8, 1.9, 2, 2, 1.9, 1.5, 1.9, 1.8, 1.7, 1.8, 1.9, 1.9, 2

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Move the active selection to E20, matching the saved selection state
$ws.Range("E20").Select()
